$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-11 with new values (columns B-K)
$ws.Cells.Item(2, 2).Value = 0.28918045258416036
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0

$ws.Cells.Item(3, 2).Value = 0.29157142997919322
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0.00018342124205140566
$ws.Cells.Item(3, 5).Value = 0.000036376800139662587
$ws.Cells.Item(3, 6).Value = -0.000076067278141884526
$ws.Cells.Item(3, 7).Value = 0.000038348851694202055
$ws.Cells.Item(3, 8).Value = -0.000015015021977388299
$ws.Cells.Item(3, 9).Value = 0.00018293432122544763
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = -0.00040277925099047041

$ws.Cells.Item(4, 2).Value = 0.28743842958423926
$ws.Cells.Item(4, 3).Value = -0.00032848497608088916
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = -0.0000075417958425524894
$ws.Cells.Item(4, 6).Value = 0.0000012179357435250763
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = -0.00003372976759955384
$ws.Cells.Item(4, 9).Value = -0.00055406131257532804
$ws.Cells.Item(4, 10).Value = -0.00004559510396636705
$ws.Cells.Item(4, 11).Value = 0.000059777127791371676

$ws.Cells.Item(5, 2).Value = 0.29861280266031204
$ws.Cells.Item(5, 3).Value = 0.0012256390847313974
$ws.Cells.Item(5, 4).Value = 0.000049722471844633673
$ws.Cells.Item(5, 5).Value = -0.000071073350177880035
$ws.Cells.Item(5, 6).Value = 0.00034930460755964288
$ws.Cells.Item(5, 7).Value = -0.00061073732516629751
$ws.Cells.Item(5, 8).Value = -0.000014910591041383004
$ws.Cells.Item(5, 9).Value = 0.0000096753796497450772
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0.000066252184312498574

$ws.Cells.Item(6, 2).Value = 0.28896203069513654
$ws.Cells.Item(6, 3).Value = 0.0056912862460800039
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0.000041505710935509671
$ws.Cells.Item(6, 6).Value = -0.000058578573851324743
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = -0.00004697425129123681
$ws.Cells.Item(6, 9).Value = -0.001671452166952321
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = -0.00000075493792067726773

$ws.Cells.Item(7, 2).Value = 0.28196001969777795
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = -0.0015496371046907684
$ws.Cells.Item(7, 5).Value = -0.00013656846288508028
$ws.Cells.Item(7, 6).Value = -0.00072841770695212789
$ws.Cells.Item(7, 7).Value = 0.00019974844129421287
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.000076632365911690766
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = -0.00023105715458537901

$ws.Cells.Item(8, 2).Value = 0.27874960665514231
$ws.Cells.Item(8, 3).Value = 0.0042925982173278553
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = -0.000027093959137777708
$ws.Cells.Item(8, 6).Value = -0.00022216635912498519
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = -0.000010895124556589706
$ws.Cells.Item(8, 9).Value = -0.00062436326470897909
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0.00016813559023215463

$ws.Cells.Item(9, 2).Value = 0.26311781192127431
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0.00031595546212201934
$ws.Cells.Item(9, 5).Value = -0.00093868145605056832
$ws.Cells.Item(9, 6).Value = -0.0055899936746864271
$ws.Cells.Item(9, 7).Value = 0.00007456828516663226
$ws.Cells.Item(9, 8).Value = -0.00016275897599123392
$ws.Cells.Item(9, 9).Value = 0.000012657203450990205
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = -0.00015011602409586988

$ws.Cells.Item(10, 2).Value = 0.28321636405378825
$ws.Cells.Item(10, 3).Value = 0.0077001772948936996
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = -0.0000043493005985976595
$ws.Cells.Item(10, 6).Value = -0.00015415012392969144
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = -0.0000026443349934821676
$ws.Cells.Item(10, 9).Value = 0.00026566884120989931
$ws.Cells.Item(10, 10).Value = -0.0017601302247825504
$ws.Cells.Item(10, 11).Value = -0.000002896520223005794

$ws.Cells.Item(11, 2).Value = 0.28913798681628533
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.0036507697710153292
$ws.Cells.Item(11, 5).Value = -0.00047364965985033852
$ws.Cells.Item(11, 6).Value = -0.0021468439024233789
$ws.Cells.Item(11, 7).Value = 0.00021072328399143937
$ws.Cells.Item(11, 8).Value = -0.00013069637814786637
$ws.Cells.Item(11, 9).Value = 0.00045386490399289346
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = -0.0026395142134688077

# Add new row 12 for 2025-08-30
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2025-08-30"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).Value = 0.31040215093789242
$ws.Cells.Item(12, 3).Value = -0.011990510800249612
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = -0.0000074353002106160982
$ws.Cells.Item(12, 6).Value = 0.0000021018870078549092
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = -0.00000034412235264679504
$ws.Cells.Item(12, 9).Value = -0.0019232344723581691
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0.000082187596686822406
